$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 534 ("「クネイトラ県、蜂の巣を約60%失う」" post), shifting all rows below it up by one.
$ws.Rows.Item(534).Delete()
